$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the placeholder/anonymized student names with real ones.
$ws.Range("B2").Value = "Darshan"
$ws.Range("B3").Value = "Saman"
$ws.Range("B4").Value = "Shail"
$ws.Range("B5").Value = "Vaibhav"
$ws.Range("B6").Value = "John"

# Leave the selection where the author left it.
$ws.Range("C6").Select()
